$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("raw_data_1")

$data = New-Object 'object[,]' 16,11
$data[0,0] = 2019
$data[0,1] = 8
$data[0,2] = 26
$data[0,3] = "LC"
$data[0,4] = "I"
$data[0,5] = 19
$data[0,6] = 1
$data[0,7] = 0.86
$data[0,8] = 16.92
$data[0,9] = 0.86
$data[0,10] = 1.06
$data[1,0] = 2019
$data[1,1] = 8
$data[1,2] = 26
$data[1,3] = "LC"
$data[1,4] = "I"
$data[1,5] = 19
$data[1,6] = 2
$data[1,7] = 0.86
$data[1,8] = 13.7
$data[1,9] = 0.86
$data[1,10] = 1.08
$data[2,0] = 2019
$data[2,1] = 8
$data[2,2] = 26
$data[2,3] = "LC"
$data[2,4] = "I"
$data[2,5] = 19
$data[2,6] = 3
$data[2,7] = 0.86
$data[2,8] = 10.94
$data[2,9] = 0.86
$data[2,10] = 1.06
$data[3,0] = 2019
$data[3,1] = 8
$data[3,2] = 26
$data[3,3] = "LC"
$data[3,4] = "I"
$data[3,5] = 19
$data[3,6] = 4
$data[3,7] = 0.86
$data[3,8] = 15.02
$data[3,9] = 0.86
$data[3,10] = 1.86
$data[4,0] = 2019
$data[4,1] = 8
$data[4,2] = 26
$data[4,3] = "LC"
$data[4,4] = "I"
$data[4,5] = 2
$data[4,6] = 2
$data[4,7] = 0.86
$data[4,8] = 20.04
$data[4,9] = 1.84
$data[4,10] = 1.22
$data[5,0] = 2019
$data[5,1] = 8
$data[5,2] = 26
$data[5,3] = "LC"
$data[5,4] = "I"
$data[5,5] = 2
$data[5,6] = 4
$data[5,7] = 0.86
$data[5,8] = 17.74
$data[5,9] = 2.86
$data[5,10] = 1.6
$data[6,0] = 2019
$data[6,1] = 8
$data[6,2] = 26
$data[6,3] = "LC"
$data[6,4] = "I"
$data[6,5] = 2
$data[6,6] = 3
$data[6,7] = 0.86
$data[6,8] = 16.46
$data[6,9] = 2.26
$data[6,10] = 0.96
$data[7,0] = 2019
$data[7,1] = 8
$data[7,2] = 26
$data[7,3] = "LC"
$data[7,4] = "I"
$data[7,5] = 2
$data[7,6] = 1
$data[7,7] = 0.86
$data[7,8] = 12.02
$data[7,9] = 1.36
$data[7,10] = 0.86
$data[8,0] = 2019
$data[8,1] = 8
$data[8,2] = 26
$data[8,3] = "LC"
$data[8,4] = "I"
$data[8,5] = 1
$data[8,6] = 4
$data[8,7] = 0.86
$data[8,8] = 21.28
$data[8,9] = 4.53
$data[8,10] = 2.15
$data[9,0] = 2019
$data[9,1] = 8
$data[9,2] = 26
$data[9,3] = "LC"
$data[9,4] = "I"
$data[9,5] = 1
$data[9,6] = 2
$data[9,7] = 0.86
$data[9,8] = 17.69
$data[9,9] = 4.99
$data[9,10] = 2.49
$data[10,0] = 2019
$data[10,1] = 8
$data[10,2] = 26
$data[10,3] = "LC"
$data[10,4] = "I"
$data[10,5] = 1
$data[10,6] = 1
$data[10,7] = 0.86
$data[10,8] = 23.13
$data[10,9] = 2.38
$data[10,10] = 7.03
$data[11,0] = 2019
$data[11,1] = 8
$data[11,2] = 26
$data[11,3] = "LC"
$data[11,4] = "I"
$data[11,5] = 1
$data[11,6] = 3
$data[11,7] = 0.86
$data[11,8] = -999
$data[11,9] = 0.86
$data[11,10] = 1.25
$data[12,0] = 2019
$data[12,1] = 8
$data[12,2] = 26
$data[12,3] = "LC"
$data[12,4] = "I"
$data[12,5] = 12
$data[12,6] = 1
$data[12,7] = 0.86
$data[12,8] = -999
$data[12,9] = 2.49
$data[12,10] = 2.15
$data[13,0] = 2019
$data[13,1] = 8
$data[13,2] = 26
$data[13,3] = "LC"
$data[13,4] = "I"
$data[13,5] = 12
$data[13,6] = 4
$data[13,7] = 0.86
$data[13,8] = -999
$data[13,9] = 1.59
$data[13,10] = 0.86
$data[14,0] = 2019
$data[14,1] = 8
$data[14,2] = 26
$data[14,3] = "LC"
$data[14,4] = "I"
$data[14,5] = 12
$data[14,6] = 3
$data[14,7] = 0.86
$data[14,8] = -999
$data[14,9] = 0.86
$data[14,10] = 0.86
$data[15,0] = 2019
$data[15,1] = 8
$data[15,2] = 26
$data[15,3] = "LC"
$data[15,4] = "I"
$data[15,5] = 12
$data[15,6] = 2
$data[15,7] = 0.86
$data[15,8] = -999
$data[15,9] = 4.08
$data[15,10] = 0.86

$ws.Range("A2:K17").Value = $data

$ws.Activate()
[void]$ws.Range("O20").Select()
